$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Average_Lifetime")

# Modified dummy lifetime value to 2.5 years to make the graphs make more sense
$ws.Range("B2").Value = 2.5

# Move/update the active selection on the sheet to match the saved view state
$ws.Range("B3").Select()
